# Apply the edits described in the diff:
# 1. Update cell F5's value from 1706219962 to 1706239962 (a last_update
#    unix-timestamp bump for the "pdf_online"/row-5 record).
# 2. Update the sheet view: scroll so column C becomes the left-most
#    visible column (topLeftCell="C1") and move the active selection from
#    A9 to E6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the F5 cell value
$ws.Range("F5").Value = 1706239962

# 2. Update the view state: scroll the window and change the selection
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 3   # column C is the 3rd column -> topLeftCell="C1"
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E6").Select()
